$wb = $excel.ActiveWorkbook

# --- Sheet: Summary ---
$ws1 = $wb.Worksheets.Item("Summary")
$ws1.Range("B2").Value = 0.3558718861209965
$ws1.Range("C2").Value = 0.06958762886597938
$ws1.Range("D2").Value = 0.9642857142857143
$ws1.Range("E2").Value = 0.1298076923076923
$ws1.Range("F2").Value = 0.27
$ws1.Range("G2").Value = 0.6452205882352942
$ws1.Range("H2").Value = 0.795211342964152
$ws1.Range("I2").Value = 27
$ws1.Range("J2").Value = 361
$ws1.Range("K2").Value = 173
$ws1.Range("L2").Value = 1

# --- Sheet: Classification Report ---
$ws2 = $wb.Worksheets.Item("Classification Report")

$ws2.Range("B2").Value = 0.9942528735632183
$ws2.Range("C2").Value = 0.3239700374531835
$ws2.Range("D2").Value = 0.4887005649717514

$ws2.Range("B3").Value = 0.06958762886597938
$ws2.Range("C3").Value = 0.9642857142857143
$ws2.Range("D3").Value = 0.1298076923076923

$ws2.Range("B4").Value = 0.3558718861209965
$ws2.Range("C4").Value = 0.3558718861209965
$ws2.Range("D4").Value = 0.3558718861209965
$ws2.Range("E4").Value = 0.3558718861209965

$ws2.Range("B5").Value = 0.5319202512145988
$ws2.Range("C5").Value = 0.6441278758694489
$ws2.Range("D5").Value = 0.3092541286397218

$ws2.Range("B6").Value = 0.9481841425106868
$ws2.Range("C6").Value = 0.3558718861209965
$ws2.Range("D6").Value = 0.4708197812803036

# --- Sheet: Confusion Matrix ---
$ws3 = $wb.Worksheets.Item("Confusion Matrix")

$ws3.Range("B2").Value = 173
$ws3.Range("C2").Value = 361

$ws3.Range("B3").Value = 1
$ws3.Range("C3").Value = 27
